$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 12: C12 ("cancel" row) JA column used to duplicate the EN string;
# give it its own proper Japanese translation "戻る".
# ---------------------------------------------------------------------------
$ws.Range("C12").Value = "戻る"

# ---------------------------------------------------------------------------
# Rows 13-19 were empty placeholder rows. Fill them with new translation
# keys/values and restyle them to match the rest of the table.
# ---------------------------------------------------------------------------

# Row 13: error_header / Error
$ws.Range("A4").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("A13").Value = "error_header"
$ws.Range("B13").Value = "Error"

# Row 14: connection_error / Connection error
$ws.Range("A4").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("A14").Value = "connection_error"
$ws.Range("B14").Value = "Connection error"

# Row 15: failed_header / Failed  (uses the alternate border variant)
$ws.Range("A4").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A15").Borders.Color = 0xa5a5a5
$ws.Range("A15").Borders.Item(7).Color = 0x3f3f3f
$ws.Range("C4").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("A15").Value = "failed_header"
$ws.Range("B15").Value = "Failed"

# Row 16: incorrect_password_error / Incorrect password (same variant as row 15)
# Reuse the style already built for row 15 instead of re-deriving the border
# colours, so no extra transient style slots are produced.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("B15").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("A16").Value = "incorrect_password_error"
$ws.Range("B16").Value = "Incorrect password"

# Row 17: recover_password_alert_header / Nice (taller row)
$ws.Range("A4").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("A17").Value = "recover_password_alert_header"
$ws.Range("B17").Value = "Nice"
$ws.Rows.Item(17).RowHeight = 27

# Row 18: recover_password_alert_msg / Sent (taller row)
$ws.Range("A4").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("A18").Value = "recover_password_alert_msg"
$ws.Range("B18").Value = "Sent"
$ws.Rows.Item(18).RowHeight = 27

# Row 19: alert_button_ok / OK
$ws.Range("A4").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("A19").Value = "alert_button_ok"
$ws.Range("B19").Value = "OK"
